$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 145, shifting existing rows 145:171 down to 146:172.
$ws.Rows.Item(145).Insert()

# Populate the newly inserted row 145 with the new record.
$ws.Cells.Item(145, 1).Value = 11
$ws.Cells.Item(145, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(145, 3).Value = "Bíobío"
$ws.Cells.Item(145, 4).Value = 44946
$ws.Cells.Item(145, 4).Style = $ws.Cells.Item(144, 4).Style
$ws.Cells.Item(145, 4).NumberFormat = $ws.Cells.Item(144, 4).NumberFormat
$ws.Cells.Item(145, 5).Value = 8
$ws.Cells.Item(145, 6).Value = "Fruta"
$ws.Cells.Item(145, 7).Value = 100109
$ws.Cells.Item(145, 8).Value = "Uva"
$ws.Cells.Item(145, 9).Value = 100109001
$ws.Cells.Item(145, 10).Value = "Uva"
$ws.Cells.Item(145, 11).Value = "Superior Seedless"
$ws.Cells.Item(145, 12).Value = "Primera"
$ws.Cells.Item(145, 13).Value = 200
$ws.Cells.Item(145, 14).Value = 10000
$ws.Cells.Item(145, 15).Value = 11000
$ws.Cells.Item(145, 16).Value = 10500
$ws.Cells.Item(145, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(145, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(145, 19).Value = 1050
$ws.Cells.Item(145, 20).Value = 10
